$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column before column D. This shifts the existing "Enum" (D)
# and "Bool" (E) columns one to the right (to E and F), carrying their
# widths/styles with them automatically.
# ---------------------------------------------------------------------------
$ws.Range("D1").EntireColumn.Insert()

# New column header ("Date"), inserted between "Date Time" and "Enum".
$ws.Range("D1").Value = "Date"

# The "Date Time" example value now includes a time-of-day component.
$ws.Range("C3").Value = "01/01/2020 10:15:05"

# New "Date" example value. Copy the formatting (alignment/wrap/number
# format) from the neighbouring cell first so the new cell's style matches
# the rest of the data row, then enter the literal as a text formula and
# immediately convert it down to a plain value so Excel cannot silently
# reinterpret the "01/01/2020" literal as a date serial number.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Formula = '="01/01/2020"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column widths: the new column gets its own width; the columns that were
# shifted right already retained their original widths via the insert.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 8.45

# ---------------------------------------------------------------------------
# Re-apply the AutoFilter over the new, wider range (A1:F3) and refresh the
# hidden _FilterDatabase defined name that Excel keeps in sync with it.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:F3").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$3"
    }
}
